{"js": "// Add five new paragraphs right after the paragraph that ends the\n// \"*: Relationship: Aggregated Relation Statement Relation (Object)\n// Roles / Kinds.\" note (and right before the following \"Model:\" block):\n//\n//   1. (blank paragraph)\n//   2. \"Intension / Extension: S / O.\"\n//   3. \"sub / super hiers, containment: P / O.\"\n//   4. (blank paragraph)\n//   5. \"Relation reification: Relation statement object: relation\n//       instance. a: Role / b: Kind: relation ends. (a): Role reifying\n//       rel attrs / values. (b): Kind Resource reifying rel subject\n//       (rel players).\"\n\nconst anchorText = \"*: Relationship: Aggregated Relation Statement Relation (Object) Roles / Kinds.\";\n\n// Locate the anchor paragraph by searching for a distinctive slice of its\n// text (search() has a length cap, so use a safe, unique substring).\nconst searchResults = context.document.body.search(\"Aggregated Relation Statement\", { matchCase: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Anchor paragraph not found: \" + anchorText);\n}\n\nlet anchor = searchResults.items[0].paragraphs.getFirst();\n\nconst newParagraphTexts = [\n  \"\",\n  \"Intension / Extension: S / O.\",\n  \"sub / super hiers, containment: P / O.\",\n  \"\",\n  \"Relation reification: Relation statement object: relation instance. a: Role / b: Kind: relation ends. (a): Role reifying rel attrs / values. (b): Kind Resource reifying rel subject (rel players).\"\n];\n\nfor (const text of newParagraphTexts) {\n  anchor = anchor.insertParagraph(text, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Add five new paragraphs right after the paragraph that ends the\n# \"*: Relationship: Aggregated Relation Statement Relation (Object)\n# Roles / Kinds.\" note (and right before the following \"Model:\" block):\n#\n#   1. (blank paragraph)\n#   2. \"Intension / Extension: S / O.\"\n#   3. \"sub / super hiers, containment: P / O.\"\n#   4. (blank paragraph)\n#   5. \"Relation reification: Relation statement object: relation\n#       instance. a: Role / b: Kind: relation ends. (a): Role reifying\n#       rel attrs / values. (b): Kind Resource reifying rel subject\n#       (rel players).\"\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph robustly via Find rather than a hard-coded\n# paragraph index.\n$searchRange = $d.Content\n$searchRange.Find.Text = \"Aggregated Relation Statement\"\n$found = $searchRange.Find.Execute()\nif (-not $found) {\n    throw \"Anchor paragraph not found\"\n}\n$anchor = $searchRange.Paragraphs(1)\n\n$newParagraphTexts = @(\n    \"\",\n    \"Intension / Extension: S / O.\",\n    \"sub / super hiers, containment: P / O.\",\n    \"\",\n    \"Relation reification: Relation statement object: relation instance. a: Role / b: Kind: relation ends. (a): Role reifying rel attrs / values. (b): Kind Resource reifying rel subject (rel players).\"\n)\n\nforeach ($text in $newParagraphTexts) {\n    $anchor.Range.InsertParagraphAfter()\n    $anchor = $anchor.Next()\n    if ($text -ne \"\") {\n        $anchor.Range.Text = $text\n    }\n}\n"}
